# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" on the Overview sheet and the
# "Latest Handoff Datetime" on the zh-cn sheet to the new handoff timestamps,
# and mark the "Ready for handoff" rows' Priority as "ht" (handoff type) on
# both the zh-cn and de-de target-language sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for rows 9-14
for ($row = 9; $row -le 14; $row++) {
    $overview.Range("G$row").Value = "2016-08-29 08:23:11"
}

# zh-cn sheet: "Latest Handoff Datetime" (column H) for rows 9-14
for ($row = 9; $row -le 14; $row++) {
    $zhcn.Range("H$row").Value = "2016-08-29 08:23:01"
}

# zh-cn and de-de sheets: set Priority (column E) to "ht" for the rows that
# are "Ready for handoff" (rows 9-14)
for ($row = 9; $row -le 14; $row++) {
    $zhcn.Range("E$row").Value = "ht"
    $dede.Range("E$row").Value = "ht"
}
